$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update prices in D29:D31
$ws.Range("D29").Value = 598
$ws.Range("D30").Value = 640
$ws.Range("D31").Value = 815
